$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Asterix - 1a (nueva lavada)" dated
# 2021-08-05 (serial 44413). It belongs right after the existing row for
# 2021-01-20 (row 17) and before the 2021-01-26 record, so insert a fresh
# row at position 18 - this pushes every following record down by one row.
$ws.Rows("18").Insert()

# Populate the newly inserted row 18 with the new record's data.
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C18").Value = "Arica y Parinacota"
$ws.Range("D18").Value = 44413
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 100114001
$ws.Range("G18").Value = "Papa"
$ws.Range("H18").Value = "Asterix"
$ws.Range("I18").Value = "1a (nueva lavada)"
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 11000
$ws.Range("L18").Value = 12000
$ws.Range("M18").Value = 11500
$ws.Range("N18").Value = "$/saco 25 kilos"
$ws.Range("O18").Value = "Provincia de Melipilla"
$ws.Range("P18").Value = 460
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
